$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D-column price cells as Text first so Excel does not auto-convert
# numeric-looking strings (e.g. "1.003", "29.019.83") into Double values,
# then restore the default "Normal" style so no stray formatting is left behind.

$dCells = @("D2","D3","D4","D5","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.019.83"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "1.910.18"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "325.31"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "0.4600"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value = "0.3871"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("D9").Value = "0.07815"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").Value = "0.9871"
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("D11").Value = "21.86"
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").Value = "1.894.37"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").Value = "5.755"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").Value = "7.002"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "0.07035"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "87.34"
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "0.000009909"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").Value = "17.01"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "0.9995"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "29.039.23"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").Value = "5.350"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").Value = "11.08"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "2.140.56"
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("D25").Value = "2.084"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").Value = "155.92"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").Value = "19.36"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").Value = "5.857"
$ws.Range("E28").Value = "  -3.09%  "
$ws.Range("D29").Value = "118.03"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "1.854"
$ws.Range("E30").Value = "  -4.81%  "
$ws.Range("D31").Value = "0.09298"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").Value = "0.8767"
$ws.Range("E32").Value = "  -4.61%  "
$ws.Range("D33").Value = "5.177"
$ws.Range("E33").Value = "  -3.05%  "
$ws.Range("D34").Value = "1.310"
$ws.Range("E34").Value = "  -2.91%  "
$ws.Range("D35").Value = "3.131"
$ws.Range("E35").Value = "  -4.58%  "
$ws.Range("D36").Value = "0.05760"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").Value = "1.171"
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("D38").Value = "0.02084"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").Value = "1.000"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "0.5678"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").Value = "7.632"
$ws.Range("E41").Value = "  -3.04%  "
$ws.Range("D42").Value = "0.1805"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").Value = "9.686"
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("D46").Value = "0.5316"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").Value = "2.177"
$ws.Range("E47").Value = "  -5.17%  "
$ws.Range("D48").Value = "0.06923"
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("D49").Value = "1.831"
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("D50").Value = "2.556"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "112.43"
$ws.Range("E51").Value = "  -0.53%  "

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# Row 44/45: PEPE and EnergySwap swap places (A/rank stays the same)
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"

$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.000002897"
$ws.Range("E44").Value = "  +89.31%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "11.80"
$ws.Range("E45").Value = "  -1.65%  "

$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
